# [ADDITIONAL SCRAPING] add a "Player Info" sheet and replace the
# MATCH_CARD_LINK columns (full scorecard URLs) with a compact
# MATCH_CODE column (just the numeric code) on both stats sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert a new "Player Info" sheet in front of "ODI Batting" and
#    fill in the player's basic info.
# ---------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$playerInfo.Range("A2").Value = "'4650"
$playerInfo.Range("B2").Value = "Jhye Avon Richardson"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# ---------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (column D) -> MATCH_CODE, values
#    collapsed from the full scorecard URL down to just the numeric
#    MatchCode query parameter.
#    (re-fetch the sheet by name -- inserting a sheet shifts what a
#    previously-grabbed worksheet reference now points at)
# ---------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$lastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Text
    if ([string]::IsNullOrEmpty($link)) { continue }
    if ($link -match 'MatchCode=(\d+)') {
        $cell.Value = "'" + $matches[1]
    }
}

# ---------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (column B) -> MATCH_CODE, same
#    URL -> numeric-code collapse.
# ---------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$lastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Text
    if ([string]::IsNullOrEmpty($link)) { continue }
    if ($link -match 'MatchCode=(\d+)') {
        $cell.Value = "'" + $matches[1]
    }
}

$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1").Select() | Out-Null
